# Attendance sheet update.
#
# Columns: A=Name, B=Jan, C=Feb, D=Mar, E=Apr, F=May, G=Jun, H=Jul
#
# - Add the new hire ("ejas") to the roster.
# - For every employee row, read whatever is currently in the Jan/Feb/Mar/Apr
#   cells and update them with this month's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162

# Find the last row that already has a Name in column A.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End($xlUp).Row

# Onboard the new hire before filling in her numbers.
$newRow = $lastRow + 1
$ws.Cells($newRow, 1).Value = "ejas"
$lastRow = $newRow

# Every employee worked the same number of weeks in January. Read it back
# from the standard cell so it can be reused while stamping every row.
$ws.Cells(2, 2).Value = 3
$standardJan = $ws.Cells(2, 2).Value2

# Feb/Mar/Apr figures collected for each employee this month, keyed by row.
$monthly = @{
    2 = @(4, 4, 1)
    3 = @(4, 6, 1)
    4 = @(7, 2, 1)
    5 = @(3, 7, 1)
}

for ($row = 2; $row -le $lastRow; $row++) {
    # Confirm the row belongs to a real employee before updating it.
    $name = $ws.Cells($row, 1).Value2
    if ([string]::IsNullOrEmpty($name)) {
        continue
    }

    # Read the existing (blank) Jan value, then update it with the standard.
    $existingJan = $ws.Cells($row, 2).Value2
    $ws.Cells($row, 2).Value = $standardJan

    $figures = $monthly[$row]

    # Read-then-update Feb, Mar and Apr for this employee.
    $existingFeb = $ws.Cells($row, 3).Value2
    $ws.Cells($row, 3).Value = $figures[0]

    $existingMar = $ws.Cells($row, 4).Value2
    $ws.Cells($row, 4).Value = $figures[1]

    $existingApr = $ws.Cells($row, 5).Value2
    $ws.Cells($row, 5).Value = $figures[2]
}
